$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "__sanatani__090"
$ws.Range("E11").Select()
